# data-quality-report: Add for list fields.
# Renames the existing "single value" fields sheet and adds a sibling sheet
# that holds the download/report configuration for "list" (multi-value)
# fields, mirroring the layout of the first sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Sheet1" -> "Fields with Single value" -------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Fields with Single value"

# --- Sheet 2: duplicate sheet 1 (keeps page setup/margins/fonts/column
# widths/cell styles) so we only need to swap in the list-field content,
# then rename and reposition it right after sheet 1. ---------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Fields with Multiple Values"

# Push the copied (already-styled) header/data rows down from 7-8 to 9-10,
# so the new sheet's content lines up the same way the original report
# template does it, while keeping each cell's original formatting.
$ws2.Rows.Item(7).Insert() | Out-Null
$ws2.Rows.Item(7).Insert() | Out-Null

# Column widths for the new sheet's two columns.
$ws2.Columns.Item(1).ColumnWidth = 40.57
$ws2.Columns.Item(2).ColumnWidth = 34.32

# Header row (bold) + data/placeholder row, same pattern as sheet 1.
$ws2.Range("A9").Value = "Field"
$ws2.Range("B9").Value = "Projects with at least one data item"

$ws2.Range("A10").Value = "SPREADSHEETFORM:DOWN:list_fields:field/title"
$ws2.Range("B10").Value = "SPREADSHEETFORM:DOWN:list_fields:count_public_projects_with_at_least_one_public_value"

# Reset selection on the new sheet, then land back on sheet 1 (A1) which
# stays the active tab/selection.
$ws2.Range("A1").Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
